$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 56 data rows (rows 2-57). The updated data dictionary
# has 60 data rows (rows 2-61), so first grow the sheet by inserting 4 blank rows
# right after the last existing data row.
$ws.Rows("58:61").Insert()

# Final contents for columns A (index), B (name), C (label), D (valueType)
# for each data row 2-61.
$data = @(
  @(1, "id", "ID", "integer"),
  @(2, "Geschlecht", "Sex", "integer"),
  @(3, "Alter_BE", "Age at blood sampling", "decimal"),
  @(4, "IPAQ_MET_Total", "IPAQ_MET_Total", "decimal"),
  @(5, "BMI", "Body Mass Index (BMI)", "decimal"),
  @(6, "kcal_zb_tag_NCI", "Energy – incl. energy from dietary fiber", "decimal"),
  @(7, "Kart_sum_NCI", "Potatoes and potato products (incl. Fresh potatoes, heated and processed products like French fries, potato pancakes, potato chips)", "decimal"),
  @(8, "Gem02_1_NCI", "Vegetables and vegetable products", "decimal"),
  @(9, "Gem02_2_NCI", "Mushrooms", "decimal"),
  @(10, "Gem02_Hu_NCI", "Legumes/Pulses (incl. lentils, chickpeas, green peas, white/broad beans, kidney/soy beans; exception: green beans; incl. canned pulses and as sauce)", "decimal"),
  @(11, "Obst_sum_NCI", "Fruit and fruit products (without juice) including unsweetened frozen fruit, fruit sauces and processed products like sweetened or heated fruit, canned fruit, dried fruit", "decimal"),
  @(12, "Nusa_sum_NCI", "Nuts and seeds (Hazelnuts, peanuts, almonds etc. or processed products like peanut butter, salted or roasted nuts)", "decimal"),
  @(13, "Milc_sum_NCI", "Total Dairy incl. milk and dairy products, cheese and curd cheese (incl. variants from goat, sheep)", "decimal"),
  @(14, "milk_NCI", "Milk", "decimal"),
  @(15, "milkbased_bev_NCI", "milk-based beverages incl. cocoa drinks, milkshakes", "decimal"),
  @(16, "Fermented_milk_NCI", "Fermented milk products", "decimal"),
  @(17, "quark_curd_NCI", "quark and curd cheese", "decimal"),
  @(18, "cheeses_NCI", "Cheese", "decimal"),
  @(19, "Sues13_6_NCI", "Creams and desserts including pudding, semolina, tiramisu", "decimal"),
  @(20, "Brot_sum_NCI", "Sum of whole grain bread/rolls, White, brown or multigrain bread/rolls, Cereals and cereal products, Pastries", "decimal"),
  @(21, "Flours_milled_prod_NCI", "Flours, milled products", "decimal"),
  @(22, "Pasta_rice_NCI", "Pasta, rice", "decimal"),
  @(23, "Brot01_1_NCI", "Sum of Wholemeal bread/rolls, White, brown or multigrain bread/rolls (incl. Crisp bread)", "decimal"),
  @(24, "Bread_NCI", "Bread", "decimal"),
  @(25, "Crispbread_NCI", "Crispbread", "decimal"),
  @(26, "Breakfast_cereals_NCI", "Breakfast cereals", "decimal"),
  @(27, "Crackers_NCI", "Crackers, breadsticks", "decimal"),
  @(28, "Fine_bakery_NCI", "Fine bakery wares", "decimal"),
  @(29, "Flei_sum_NCI", "Sum of meat and meat products unsmoked (including roast, goulash, schnitzel, beef olive, minced meat, meat balls, ground meat sauce) & sausages and meat products smoked (including salami, liver sausage, ham, bacon, cured pork, meat loaf, bratwurst)", "decimal"),
  @(30, "Flei_Wurst_NCI", "Sausages and meat products smoked including salami, liver sausage, ham, bacon, cured pork, meat loaf, bratwurst", "decimal"),
  @(31, "Fish_sum_NCI", "Fish and fish products including salt- and freshwater fish, shrimps, mussels, snails, processed products like caviar, tinned fish", "decimal"),
  @(32, "Eier_sum_NCI", "Eggs and egg products including scrambled, fried and boiled eggs, omelets, soufflé", "decimal"),
  @(33, "Fett_sum_NCI", "fats & oils (excl: in salad dressings) - incl. Animal-based fats, plant-based fats and not specified fats/oils", "decimal"),
  @(34, "veg_fatsoils_NCI", "Vegetable fats and oils", "decimal"),
  @(35, "Butter_NCI", "Butter", "decimal"),
  @(36, "Marg_Streichf_NCI", "Margarine and similar", "decimal"),
  @(37, "Animal_fatsoils_NCI", "Animal fats excluding butter, fish oil", "decimal"),
  @(38, "Sues_sum_NCI", "Sweets (sum)", "decimal"),
  @(39, "chocolate_sweets_NCI", "Sweets including chocolate, sweets with chocolate coating, confectioneries with chocolate", "decimal"),
  @(40, "nonchoc_sweets_NCI", "Sweets without chocolate such as  candies, fruit gums, cereal bars", "decimal"),
  @(41, "Sues13_2_NCI", "Ice cream", "decimal"),
  @(42, "Cakes_NCI", "Cakes", "decimal"),
  @(43, "Getr15_1_NCI", "non-alcoholic beverages (sum)", "decimal"),
  @(44, "Getr15_14_NCI", "Fruit juices/nectars (incl. Spritzer) - excl. vegetable juices", "decimal"),
  @(45, "Getr15_16_NCI", "Soft drinks e.g. lemonades, bitter lemon (incl. low-calorie drinks)", "decimal"),
  @(46, "Getr15_11_NCI", "Water", "decimal"),
  @(47, "blackgreentea_NCI", "black and green tea", "decimal"),
  @(48, "Kaffee_NCI", "coffee", "decimal"),
  @(49, "Getr15_13_NCI", "Herbal/fruit tea", "decimal"),
  @(50, "Getr15_2_NCI", "Alcoholic beverages (sum)", "decimal"),
  @(51, "Getr15_22_NCI", "Wine and sparkling wine", "decimal"),
  @(52, "Getr15_21_NCI", "Beer including mixed beer drinks", "decimal"),
  @(53, "Getr15_23_MW", "spirits incl. schnaps, liqueurs", "decimal"),
  @(54, "Getr15_24_MW", "Other alcoholic beverages e.g. alcopops, cocktails", "decimal"),
  @(55, "Sose_sum_NCI", "Sauces including warm and cold sauces (e.g. ketchup, salad dressing), mustard, vinegar etc.; exceptions: fruit-, vegetable-, ground meat sauce", "decimal"),
  @(56, "Supp_sum_NCI", "Soups without stews", "decimal"),
  @(57, "Sonst_sum_NCI", "Milk substitutes (e.g. soy products (milk, yoghurt, cheese), coconut-, rice-, oat-, almond milk, cereal milk with soy), meat substitutes (e.g. tofu, tempeh, soy protein, vegetarian sausages), cereal substitutes (e.g. soy flour/flakes, lupines), sweeteners, sugar substitutes, beverage powders/ -granules (e.g. cocoa powder, lemonade powder), herbs, spices, vegetarian spreads, protein powder, yeast, miso", "decimal"),
  @(58, "Veg_products_NCI", "Vegetarian products", "decimal"),
  @(59, "Soyproducts_NCI", "Soyproducts", "decimal"),
  @(60, "Suessstoffe_MW", "table sweeteners (low/no caloric sweeteners) (excl. Stevia)", "decimal")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
